$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells G1 and H1
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Copy formatting from an existing header cell (F1) so the new headers
# share the same bold/border/centered style
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Update existing values in row 2
$ws.Range("B2").Value = 0.2379811448550953
$ws.Range("C2").Value = 0.9953486525390304
$ws.Range("D2").Value = 0.3804142708913718

# Add new values G2 and H2
$ws.Range("G2").Value = 0.1239050709499376
$ws.Range("H2").Value = 0.991
